# Updated cryptos list values (price + volume) to match latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.479.11'
$ws.Range("E2").Value = '  +2.26%  '
$ws.Range("D3").Value = '1.865.67'
$ws.Range("E3").Value = '  +2.71%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.34%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.00'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.53%  '
$ws.Range("E6").Value = '  -0.37%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4662'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3736'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.51%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07390'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8887'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.64%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07962'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +6.06%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '20.06'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.83%  '
$ws.Range("D13").Value = '1.854.26'
$ws.Range("E13").Value = '  +2.06%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.431'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.09%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.602'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.54%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '92.74'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.14%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.006'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.24%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008962'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.49%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.004'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.36%  '
$ws.Range("E20").Value = '  +3.71%  '
$ws.Range("D21").Value = '27.504.86'
$ws.Range("E21").Value = '  +2.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.182'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.58'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.85%  '
$ws.Range("D24").Value = '2.086.86'
$ws.Range("E24").Value = '  +1.82%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.04'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.97%  '
$ws.Range("E26").Value = '  +2.44%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.57'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.64%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.093'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.99%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.167'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.65%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '117.30'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.83%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08916'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.70%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.014'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.91%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7518'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.50%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.161'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.03%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.494'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.92%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.664'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +10.85%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01970'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.93%  '
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.083'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.04%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05284'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.04%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.991'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.54%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.197'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.90%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5234'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.98%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1647'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.64%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.365'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.62%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4896'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.31%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.40'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.75%  '
$ws.Range("E47").Value = '  -0.37%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.670'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.44%  '
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '103.73'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.71%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06264'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.41%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '66.03'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.47%  '
